$d = $word.ActiveDocument

# 1. Remove the "Meta description" paragraph (the second paragraph of the document).
$metaPara = $d.Paragraphs.Item(2)
$metaPara.Range.Delete() | Out-Null

# 2. Insert a new paragraph, right before the final ("Create a Feature Image
#    Prompt") paragraph, that repeats the page title in bold.
$lastIndex = $d.Paragraphs.Count
$lastPara = $d.Paragraphs.Item($lastIndex)
$lastPara.Range.InsertParagraphBefore() | Out-Null

$newIndex = $d.Paragraphs.Count - 1
$newPara = $d.Paragraphs.Item($newIndex)
$newRange = $newPara.Range.Duplicate
$newRange.MoveEnd(1, -1) | Out-Null

$newParaXml = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:r/><w:r><w:rPr><w:b/></w:rPr><w:t>Play Cleopatra''s Diamonds for Free - Review</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$newRange.InsertXML($newParaXml) | Out-Null

# 3. Replace the text of the final paragraph (formerly the image-prompt
#    text) with the meta-description copy, keeping its italic run formatting.
$finalIndex = $d.Paragraphs.Count
$finalPara = $d.Paragraphs.Item($finalIndex)
$finalRange = $finalPara.Range.Duplicate
$finalRange.MoveEnd(1, -1) | Out-Null
$finalRange.Text = "Read our review of Cleopatra's Diamonds, the online slot game developed by SWINTT, and play for free. Enjoy unique bonus features, medium volatility and 96.01% RTP."
